$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.898.72"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "2.450.72"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'524.49"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "'130.77"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("D9").Value = "'0.0979"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D11").Value = "'4.95"
$ws.Range("E11").Value = "  -3.85%  "
$ws.Range("D12").Value = "'0.322"
$ws.Range("E12").Value = "  -3.09%  "
$ws.Range("D13").Value = "2.884.08"
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("D14").Value = "57.794.34"
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").Value = "'21.81"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("D17").Value = "2.446.45"
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("D18").Value = "'10.39"
$ws.Range("E18").Value = "  -2.69%  "
$ws.Range("D19").Value = "'4.15"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").Value = "'315.85"
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("D21").Value = "'6.09"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "'65.06"
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("D24").Value = "'0.408"
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").Value = "'0.156"
$ws.Range("E26").Value = "  -2.84%  "
$ws.Range("D27").Value = "'7.24"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").Value = "'172.24"
$ws.Range("E28").Value = "  +2.89%  "
$ws.Range("D29").Value = "0.0₃0737"
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("E31").Value = "  -3.76%  "
$ws.Range("D32").Value = "'6.12"
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("D35").Value = "'17.82"
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("E36").Value = "  -6.51%  "
$ws.Range("D37").Value = "'3.82"
$ws.Range("E37").Value = "  -3.88%  "
$ws.Range("D38").Value = "'36.28"
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("D39").Value = "'1.47"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").Value = "'0.797"
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("D41").Value = "'3.42"
$ws.Range("E41").Value = "  -2.15%  "
$ws.Range("D42").Value = "'265.33"
$ws.Range("E42").Value = "  -4.96%  "
$ws.Range("D43").Value = "'0.584"
$ws.Range("E43").Value = "  -2.41%  "
$ws.Range("D44").Value = "'4.81"
$ws.Range("E44").Value = "  -5.72%  "
$ws.Range("D45").Value = "'124.62"
$ws.Range("E45").Value = "  +1.01%  "
$ws.Range("D46").Value = "'0.0929"
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").Value = "'0.0211"
$ws.Range("E48").Value = "  -1.35%  "
$ws.Range("D49").Value = "'17.06"
$ws.Range("E49").Value = "  -4.38%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'16.39"
$ws.Range("E50").Value = "  -3.30%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.707.19"
$ws.Range("E51").Value = "  -2.15%  "
